# Update the cuotas archive sheet with the latest transaction data pulled
# from the Menta API (BIND / BOTON / DEVBOTON processing run).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Transaccion / Cuotas pairs, in row order starting at row 2.
$data = @(
    @(151979169, 1),
    @(726896511, 1),
    @(550338825, 1),
    @(764228034, 1),
    @(758146980, 1),
    @(171935222, 2),
    @(20645182, 3),
    @(79418400, 1),
    @(641239853, 1),
    @(110151679, 1),
    @(749067789, 1),
    @(189153073, 1),
    @(254395074, 1)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
